$d = $word.ActiveDocument
$cr = [char]13
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 1) Append a "." run to the end of the "NPX: ..." paragraph.
# ---------------------------------------------------------------------------
$pNpx = $d.Paragraphs.Item(6)
$rNpx = $d.Range($pNpx.Range.Start, $pNpx.Range.End - 1)
$rNpx.InsertAfter(".")

# ---------------------------------------------------------------------------
# 2) Append a "." run to the end of the "NPM is for management ..." paragraph.
# ---------------------------------------------------------------------------
$pNpm = $d.Paragraphs.Item(7)
$rNpm = $d.Range($pNpm.Range.Start, $pNpm.Range.End - 1)
$rNpm.InsertAfter(".")

# ---------------------------------------------------------------------------
# 3) Insert the new paragraphs (blank line, heading, and bullet content) right
#    after that last paragraph.
# ---------------------------------------------------------------------------
$pLast = $d.Paragraphs.Item(7)
$rLast = $d.Range($pLast.Range.Start, $pLast.Range.End - 1)

$lines = @(
    "",
    "3) What is babel & Web-pack?",
    "$bullet Bable is a transpiler & also a compiler.",
    "$bullet It will take ES6 input & converts it to ES5 for our browsers to understand.",
    "$bullet It also converts JSX (JavaScript Extended) language to JS.",
    "$bullet Web-pack is basically bundling tool used by react:",
    "$bullet It will bundle all the CSS & html files & will add to Index.html file.   " + " (Using <link> tags)"
)
$newText = $cr + ($lines -join $cr)
$rLast.InsertAfter($newText)

# ---------------------------------------------------------------------------
# 4) Make the "3) What is babel & Web-pack?" paragraph bold.
# ---------------------------------------------------------------------------
$pHeading = $d.Paragraphs.Item(9)
$rHeading = $pHeading.Range
$rHeading.Bold = 1
$rHeading.BoldBi = 1
